# PS2-1 / PS2-2 workbook update
# - Adds an explanatory comment about the "J" ambiguity symbol in a merged,
#   word-wrapped cell block (C13:M15).
# - Replaces the old (mostly empty) 21x21 amino-acid grid (rows 14-35) with a
#   filled-in 10x10 grid using a simplified A-J alphabet (rows 16-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: insert two new rows at row 14 so the old header (row 14) and
#    its data rows (15-35) shift down to rows 16-37.
# ---------------------------------------------------------------------------
$ws.Rows("14:15").Insert()

# ---------------------------------------------------------------------------
# 2. Comment block: merge C13:M15, enable word wrap, and write the note.
# ---------------------------------------------------------------------------
$commentRange = $ws.Range("C13:M15")
$commentRange.Merge()
$commentRange.WrapText = $true
$ws.Range("C13").Value = "J isn't a real amino acid, but can be a symbol for ambiguity between I and L. However, J is not in BLOSUM, so I assume these aren't really meant to be genuine amino acid symbols?"

# ---------------------------------------------------------------------------
# 3. New header row (row 16): simplified alphabet A..J across D16:M16.
# ---------------------------------------------------------------------------
$header = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $header.Length; $i++) {
    $col = 4 + $i   # D = 4
    $ws.Cells.Item(16, $col).Value = $header[$i]
}

# ---------------------------------------------------------------------------
# 4. Row labels (C17:C26) + matrix values for the new 10x10 grid.
# ---------------------------------------------------------------------------
$rowLabels = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $rowLabels.Length; $i++) {
    $r = 17 + $i
    $ws.Cells.Item($r, 3).Value = $rowLabels[$i]
}

# Clear any stale values left behind in the data area after the header move
# (old rows 15-16 -> new rows 17-18 carried over stray numbers/labels).
$ws.Range("D17:M26").ClearContents()

# Matrix values: column letters map to columns D(4)..M(13) per $header above.
$ws.Range("J17").Value = 7

$ws.Range("J19").Value = 7
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 9

$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 9
$ws.Range("L20").Value = 11

$ws.Range("J21").Value = 16
$ws.Range("L21").Value = 4

$ws.Range("L22").Value = 9
$ws.Range("M22").Value = 6

$ws.Range("D23").Value = 7
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 16

$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 9

$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 4
$ws.Range("I25").Value = 9

$ws.Range("F26").Value = 9
$ws.Range("I26").Value = 6

# ---------------------------------------------------------------------------
# 5. Remove the leftover amino-acid row labels (rows 27-37, old rows 25-35)
#    that no longer apply to the simplified A-J grid, while keeping the
#    existing row-height formatting (blank filler rows).
# ---------------------------------------------------------------------------
$ws.Range("C27:C37").ClearContents()

# ---------------------------------------------------------------------------
# 6. Update the selection to match the new comment block.
# ---------------------------------------------------------------------------
$ws.Range("C13:M15").Select()
